# Add a new "Barbarian" Health Point test case row to the "Health Points" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Health Points")

# Duplicate the formatting of the existing test-case row (row 3) into the new row 4
# so the new row inherits the same cell styles (number/center alignment, italic note cell, etc).
$ws.Range("A3:L3").Copy()
$ws.Range("A4:L4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new test case values (Barbarian class test case).
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = $ws.Range("B3").Value()
$ws.Range("C4").Value = $ws.Range("C3").Value()
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = -2
$ws.Range("I4").Value = "Barbarian"
$ws.Range("J4").Value = $ws.Range("J3").Value()
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = ""

# Update the sheet selection to match the new active cell after the edit.
$ws.Activate() | Out-Null
$ws.Range("L19").Select() | Out-Null
